$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 holds the sample "ModificarCliente" record. Update the client's
# identity/address data and the test-run result columns (Estado /
# Transaccion / Fecha) to reflect the new regression run.
#
# Cells C2,D2,E2,G2,H2,J2 carry the "quote-prefix" text style (s="2" /
# quotePrefix) inherited from the original cells (DNI-like numbers and
# names stored as text). A leading "'" keeps Range.Value writes typed as
# text and preserves that style instead of resetting it.
$ws.Range("C2").Value = "'11122549"
$ws.Range("D2").Value = "'Giron"
$ws.Range("E2").Value = "'Castro"
$ws.Range("G2").Value = "'Luis Giron Castro"
$ws.Range("H2").Value = "'Casado"
$ws.Range("J2").Value = "'Av. Direccion 1"

# Estado / Transaccion / Fecha for the new run.
$ws.Range("K2").Value = "FAILED"
$ws.Range("L2").Value = ""
$ws.Range("M2").Value = "26 jun. 2023, 18:21:23"

# Leave the current selection on K2, matching the saved view state.
$ws.Range("K2").Select() | Out-Null
